$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: MSV
$ws.Range("A27").Value = "MSV"
$ws.Range("A27").Font.Bold = $true
$ws.Range("B27").Value = 448
$ws.Range("C27").Value = 502
$ws.Range("D27").Value = 609
$ws.Range("F27").Value = 397
$ws.Range("G27").Value = 98

# Row 28: BSIV
$ws.Range("A28").Value = "BSIV"
$ws.Range("A28").Font.Bold = $true
$ws.Range("B28").Value = 449
$ws.Range("C28").Value = 501
$ws.Range("D28").Value = 608
$ws.Range("F28").Value = 401
$ws.Range("G28").Value = 402

$ws.Range("G28").Select()
